# Fix issue #32 (Hello X-Arguments cut off on folio 2 of Ch. 4 "Maven")
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached "fixed" date shown by the datetimeFigureOut
#    field on the slide master and on every slide layout
#    (05.02.2018 -> 20.03.2018).
# ---------------------------------------------------------------------
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Datumsplatzhalter*") {
        $shp.TextFrame.TextRange.Text = "20.03.2018"
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Datumsplatzhalter*") {
            $shp.TextFrame.TextRange.Text = "20.03.2018"
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 1 ("Chapter 4: Maven" title) - re-typing across the
#    "Chapter "/"4: " run boundary merges those two (identically
#    formatted) runs into a single "Chapter 4: " run.
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$mergedSpan = $titleRange.Characters(1, 11)
$mergedSpan.Text = "Chapter 4: "

# ---------------------------------------------------------------------
# 3) Slide 2 - enlarge/reposition the "C:\>javac ..." console output
#    box ("Rectangle 4") so the X-Arguments output is no longer cut
#    off.
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$consoleBox = $slide2.Shapes.Item(2)
$consoleBox.Left = 135.08543307086615
$consoleBox.Top = 293.9783478566929
$consoleBox.Width = 685.0287401574803
$consoleBox.Height = 181.75779527559055
